# Add a "trigger" column (B) next to the existing "stimuli" column (A),
# mirroring each stimulus filename with a "trigger_" prefixed counterpart,
# and fix the one stimuli filename that was missing its ".wav" extension.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix A16: the Keyboard stimulus entry was missing its .wav extension ---
$ws.Range("A16").Value = "Stimuli/Keyboard - major - minor - A-G-F-E 80BPM.wav"

# --- Build the new "trigger" column values (header + 15 data rows) ---
$triggerValues = @(
    "trigger",
    "Stimuli/trigger_Albums-AnaBelen_Veneo-03.wav",
    "Stimuli/trigger_Albums-Ballroom_Classics4-11.wav",
    "Stimuli/trigger_Albums-Ballroom_Classics4-12.wav",
    "Stimuli/trigger_Albums-Ballroom_Magic-09.wav",
    "Stimuli/trigger_Albums-Cafe_Paradiso-06.wav",
    "Stimuli/trigger_Albums-Cafe_Paradiso-07.wav",
    "Stimuli/trigger_Albums-Cafe_Paradiso-09.wav",
    "Stimuli/trigger_Albums-Cafe_Paradiso-10.wav",
    "Stimuli/trigger_Albums-Cafe_Paradiso-15.wav",
    "Stimuli/trigger_Albums-Cafe_Paradiso-16.wav",
    "Stimuli/trigger_Albums-Commitments-10.wav",
    "Stimuli/trigger_Albums-Fire-01.wav",
    "Stimuli/trigger_Albums-GloriaEstefan_MiTierra-04.wav",
    "Stimuli/trigger_Albums-GloriaEstefan_MiTierra-06.wav",
    "Stimuli/trigger_Keyboard - major - minor - A-G-F-E 80BPM.wav"
)

for ($i = 0; $i -lt $triggerValues.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 2).Value = $triggerValues[$i]
}

# --- Match column B's formatting (rows 2-16) to column A's; the header
#     cell B1 is left with the default style, matching A1's neighbour cell. ---
$ws.Range("A2:A16").Copy() | Out-Null
$ws.Range("B2:B16").PasteSpecial(-4122) | Out-Null
$ws.Range("A1").Select() | Out-Null

# --- Auto-size the two columns to fit their (now longer) contents ---
$ws.Columns("A:B").AutoFit() | Out-Null

Write-Host "Added trigger column with $($triggerValues.Length) rows"
